# Applies the Jan 27 2024 18:31:52 UTC cryptos-list price/volume refresh
# (GitHub Actions scheduled update) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.812.13"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.269.16"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'304.92"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "'93.01"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.530"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'32.61"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'0.0798"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "'0.113"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "2.620.72"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'14.31"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Value = "2.275.53"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'0.785"
$ws.Range("E17").Value = "  +3.57%  "
$ws.Range("D18").Value = "41.754.85"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "'12.99"
$ws.Range("E19").Value = "  +5.94%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'67.85"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'244.24"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "'2.59"
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D27").Value = "'24.05"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "'9.62"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  -5.32%  "
$ws.Range("D30").Value = "'34.94"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").Value = "'159.70"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "'5.32"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'0.0743"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").Value = "'16.94"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.106"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'2.37"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "'3.95"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'19.81"
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("D43").Value = "2.009.44"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").Value = "'2.24"
$ws.Range("E44").Value = "  +12.15%  "
$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").Value = "'10.29"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").Value = "'2.92"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "'53.51"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").Value = "'73.21"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").Value = "'1.14"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("E51").Value = "  -1.25%  "
